# Adds the next few days of COVID-19 tracking data (2/29/2020 - 3/11/2020)
# to the "wuhan", "hubei" and "china" sheets, matching the upstream
# "Add files via upload" commit.

$wb = $excel.ActiveWorkbook

function Set-RowVals {
    # Writes a run of values starting at ($row, $startCol) on $sheet, going
    # right. $null entries are skipped (leaving any existing cell alone).
    param($sheet, $row, $startCol, $vals)
    $c = $startCol
    foreach ($v in $vals) {
        if ($v -ne $null) {
            $sheet.Cells.Item($row, $c).Value = $v
        }
        $c = $c + 1
    }
}

function Set-DateCell {
    # Writes an Excel date serial into a cell and makes sure it keeps/gets
    # a date number format (so it lines up with the existing B-column dates).
    param($sheet, $row, $col, $serial, $fmt)
    $cell = $sheet.Cells.Item($row, $col)
    $cell.Value = $serial
    $cell.NumberFormat = $fmt
}

# ---------------------------------------------------------------------
# Sheet "wuhan"
# ---------------------------------------------------------------------
$wuhan = $wb.Worksheets.Item("wuhan")

Set-RowVals $wuhan 53 3 @(565, 26, 1675, 49122, 2195, 19227)

Set-DateCell $wuhan 54 2 43891 "m/d/yy"
Set-RowVals  $wuhan 54 1 @(50)
Set-RowVals  $wuhan 54 3 @(193, 32, 1958, 49315, 2227, 21185)

Set-DateCell $wuhan 55 2 43892 "m/d/yy"
Set-RowVals  $wuhan 55 1 @(51)
Set-RowVals  $wuhan 55 3 @(111, 24, 1846, 49426, 2251, 23031)

Set-DateCell $wuhan 56 2 43893 "m/d/yy"
Set-RowVals  $wuhan 56 1 @(52)
Set-RowVals  $wuhan 56 3 @(114, 31, 1859, 49540, 2282, 24890)

Set-DateCell $wuhan 57 2 43894 "m/d/yy"
Set-RowVals  $wuhan 57 1 @(53)

$wuhan.Range("G56").Select()

Write-Output "wuhan done"
